$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @(2888, 'Asosiasi FinTech Indonesia (AFTECH)', 'Asosiasi'),
    @(2889, 'Asosiasi Fintech Syariah Indonesia (AFSI)', 'Asosiasi'),
    @(2890, 'Asosiasi Blockchain Indonesia (ABI)', 'Asosiasi'),
    @(2891, 'PT Aset Digital Berkat (Tokocrypto)', 'Pedagang Aset Keuangan Digital'),
    @(2892, 'PT Bumi Santosa Cemerlang (Pluang)', 'Pedagang Aset Keuangan Digital'),
    @(2893, 'PT CTXG Indonesia Berkarya (Mobee)', 'Pedagang Aset Keuangan Digital'),
    @(2894, 'PT Indodax Nasional Indonesia (Indodax)', 'Pedagang Aset Keuangan Digital'),
    @(2895, 'PT Aset Kripto Internasional (BTSE.ID)', 'Pedagang Aset Keuangan Digital'),
    @(2896, 'PT Cipta Koin Digital (Naga Exchange)', 'Pedagang Aset Keuangan Digital'),
    @(2897, 'PT Enkripsi Teknologi Handal (Usenobi)', 'Pedagang Aset Keuangan Digital'),
    @(2898, 'PT Kagum Teknologi Indonesia (Ajaib)', 'Pedagang Aset Keuangan Digital'),
    @(2899, 'PT Kripto Maksima Koin (Floq)', 'Pedagang Aset Keuangan Digital'),
    @(2900, 'PT Mitra Kripto Sukses (Kriptosukses)', 'Pedagang Aset Keuangan Digital'),
    @(2901, 'PT Pintu Kemana Saja (Pintu)', 'Pedagang Aset Keuangan Digital'),
    @(2902, 'PT Sentra Bitwewe Indonesia (Bitwewe)', 'Pedagang Aset Keuangan Digital'),
    @(2903, 'PT Tiga Inti Utama (Triv)', 'Pedagang Aset Keuangan Digital'),
    @(2904, 'PT Upbit Exchange Indonesia (Upbit)', 'Pedagang Aset Keuangan Digital'),
    @(2905, 'PT Rekeningku Dotcom Indonesia (Reku)', 'Pedagang Aset Keuangan Digital'),
    @(2906, 'PT Teknologi Struktur Berantai (Bitwyre)', 'Pedagang Aset Keuangan Digital'),
    @(2907, 'PT Tumbuh Bersama Nano (Nanovest)', 'Pedagang Aset Keuangan Digital'),
    @(2908, 'PT Multikripto Exchange Indonesia (Koinsayang)', 'Pedagang Aset Keuangan Digital'),
    @(2909, 'PT Samuel Kripto Indonesia (Samuel Kripto)', 'Pedagang Aset Keuangan Digital'),
    @(2910, 'PT Utama Aset Digital Indonesia (Bittime)', 'Pedagang Aset Keuangan Digital'),
    @(2911, 'PT Aset Instrumen Digital (ASTAL)', 'Pedagang Aset Keuangan Digital'),
    @(2912, 'PT Central Finansial X', 'Penyelenggara Bursa Aset Keuangan Digital'),
    @(2913, 'PT Kliring Komoditi Indonesia', 'Lembaga Kliring Penjaminan dan Penyelesaian Perdagangan Aset Keuangan Digital'),
    @(2914, 'PT Kustodian Koin Indonesia', 'Pengelola Tempat Penyimpanan Aset Keuangan Digital'),
    @(2915, 'PT Tennet Depository Indonesia', 'Pengelola Tempat Penyimpanan Aset Keuangan Digital'),
    @(2916, 'PT Trusting Social Indonesia', 'Innovative Credit Scoring'),
    @(2917, 'PT Scoring Teknologi Indonesia', 'Innovative Credit Scoring'),
    @(2918, 'PT Izi Data Indonesia', 'Innovative Credit Scoring'),
    @(2919, 'PT Provenir Data Indonesia', 'Innovative Credit Scoring'),
    @(2920, 'PT Semangat Digital Bangsa', 'Innovative Credit Scoring'),
    @(2921, 'PT Prime Analytics Indonesia', 'Innovative Credit Scoring'),
    @(2922, 'PT Bangun Percaya Sosial', 'Innovative Credit Scoring'),
    @(2923, 'PT Tongdun Technology Indonesia', 'Innovative Credit Scoring'),
    @(2924, 'PT Aiforesee Inovasi Skor', 'Innovative Credit Scoring'),
    @(2925, 'PT Cloudun Technology Indonesia', 'Innovative Credit Scoring'),
    @(2926, 'PT Finture Tech Indonesia', 'Agregasi Informasi Produk dan Layanan Jasa Keuangan'),
    @(2927, 'PT Sarana Pasar Digital', 'Agregasi Informasi Produk dan Layanan Jasa Keuangan'),
    @(2928, 'PT Ringkas Asia Technology', 'Agregasi Informasi Produk dan Layanan Jasa Keuangan'),
    @(2929, 'PT Teknologi Cerdas Finansial', 'Agregasi Informasi Produk dan Layanan Jasa Keuangan'),
    @(2930, 'PT Finetiks Inovasi Indonesia', 'Agregasi Informasi Produk dan Layanan Jasa Keuangan'),
    @(2931, 'PT Hidup Ideal Sejahtera', 'Agregasi Informasi Produk dan Layanan Jasa Keuangan'),
    @(2932, 'PT Komunal Sejahtera Indonesia', 'Agregasi Informasi Produk dan Layanan Jasa Keuangan'),
    @(2933, 'PT Lendana Digitalindo Nusantara', 'Agregasi Informasi Produk dan Layanan Jasa Keuangan'),
    @(2934, 'PT Inovasi Finansial Untuk Indonesia', 'Agregasi Informasi Produk dan Layanan Jasa Keuangan'),
    @(2935, 'PT Estrend Teknologi Digital', 'Agregasi Informasi Produk dan Layanan Jasa Keuangan'),
    @(2936, 'PT Efunding Teknologi Keuangan', 'Agregasi Informasi Produk dan Layanan Jasa Keuangan'),
    @(2937, 'PT Indoartha Perkasa Sukses', 'Agregasi Informasi Produk dan Layanan Jasa Keuangan'),
    @(2938, 'PT Nex Teknologi Digital', 'Agregasi Informasi Produk dan Layanan Jasa Keuangan'),
    @(2939, 'PT Reliance Integrasi Dunia Anda', 'Agregasi Informasi Produk dan Layanan Jasa Keuangan'),
    @(2940, 'PT Agregasi Cermat Indonesia', 'Agregasi Informasi Produk dan Layanan Jasa Keuangan'),
    @(2941, 'PT Solusi Inklusi Finansial', 'Agregasi Informasi Produk dan Layanan Jasa Keuangan'),
    @(2942, 'PT Loan Market Indo', 'Agregasi Informasi Produk dan Layanan Jasa Keuangan'),
    @(2943, 'PT Xfers Indonesia Teknologi', 'Agregasi Informasi Produk dan Layanan Jasa Keuangan'),
    @(2944, 'PT Unicorn Technology Indonesia', 'Agregasi Informasi Produk dan Layanan Jasa Keuangan'),
    @(2945, 'PT Seedflex Technologies Indonesia', 'Agregasi Informasi Produk dan Layanan Jasa Keuangan'),
)

$startRow = 2889
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $item = $rows[$i]
    $ws.Cells.Item($r, 1).Value = $item[0]
    $ws.Cells.Item($r, 2).Value = $item[1]
    $ws.Cells.Item($r, 3).Value = $item[2]
}

$ws.Range("D2887").Select()
